# Revised data mapping and added resources folder
# - Source is now a CSV export ("movies_db") instead of raw text dumps from
#   separate per-source "file"/"target_db" labels.
# - Header for column A renamed from "Source DB" to "Source File".
# - Re-point the rows to the new netflix_title.csv / "IMDb movies.csv"
#   source files and the consolidated movies_db target database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "Source File"

# --- Netflix rows (2-5): file/target_db -> csv/netflix_title.csv/movies_db
$ws.Range("A2").Value = "csv"
$ws.Range("B2").Value = "netflix_title.csv"
$ws.Range("F2").Value = "movies_db"

$ws.Range("A3").Value = "csv"
$ws.Range("B3").Value = "netflix_title.csv"
$ws.Range("F3").Value = "movies_db"

$ws.Range("A4").Value = "csv"
$ws.Range("B4").Value = "netflix_title.csv"
$ws.Range("F4").Value = "movies_db"

$ws.Range("A5").Value = "csv"
$ws.Range("B5").Value = "netflix_title.csv"
$ws.Range("F5").Value = "movies_db"

# --- IMDb rows (6-7): file/target_db -> csv/IMDb movies.csv/movies_db
$ws.Range("A6").Value = "csv"
$ws.Range("B6").Value = "IMDb movies.csv"
$ws.Range("F6").Value = "movies_db"

$ws.Range("A7").Value = "csv"
$ws.Range("B7").Value = "IMDb movies.csv"
$ws.Range("F7").Value = "movies_db"

# --- Column widths (B widened / no longer auto bestFit, F widened) -----
$ws.Columns.Item(2).ColumnWidth = 14.59
$ws.Columns.Item(6).ColumnWidth = 9.75

# --- Selection moved from F21 to J11 ------------------------------------
$ws.Range("J11").Select()
